$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "0fea9b83-0187-4b9a-a30d-7fcc468cd795_fila_8.png"
$ws.Range("B5").Value = "David Orlando Mena Valverd"
$ws.Range("A6").Value = "25ca8875-48b9-4d6f-83a1-64c6c65b79f1_fila_13.png"
$ws.Range("B6").Value = "Marlon Josue Gonzales Cano"
$ws.Range("A7").Value = "2d21abef-a261-4e7f-a2a7-27a463ef0506_fila_15.png"
$ws.Range("B7").Value = "Orlando Mauricio Guevara"
$ws.Range("A8").Value = "5b439dc6-6e18-4fc1-af77-a153f1f91b73_fila_6.png"
$ws.Range("B8").Value = "Ronier Jose Rivera"
$ws.Range("A9").Value = "84817d56-0c2f-4152-9108-ba0206098ae6_fila_9.png"
$ws.Range("B9").Value = "Roman Alfonso Grios Boza"
$ws.Range("A10").Value = "a2219dad-4ed3-4a04-8e41-3f7c3e410c44_fila_11.png"
$ws.Range("B10").Value = "Eduardo Domingo Zeledon Merca"
$ws.Range("A11").Value = "a6c522d8-c904-4bd2-9498-dd63a5af504a_fila_5.png"
$ws.Range("B11").Value = "Erick Espinoza"
$ws.Range("A12").Value = "ac574bc6-1ae8-46d2-95a3-2d91342d0985_fila_14.png"
$ws.Range("B12").Value = "Angel Isaac Alvarez Quiñonez"
$ws.Range("A13").Value = "d00088f5-e013-47c4-a206-c36bd854fe67_fila_3.png"
$ws.Range("B13").Value = "Bryan Alexander Cano"
$ws.Range("A14").Value = "d0eee7e1-571b-48c0-80dc-df11791526f7_fila_12.png"
$ws.Range("B14").Value = "José Danilo Suárez"
$ws.Range("A15").Value = "d2a27921-1ef0-416c-b85e-a1a08eab12be_fila_1.png"
$ws.Range("B15").Value = "Hotep Antonio Ruiz Lezama"
$ws.Range("A16").Value = "d5fba4c7-5088-4065-9ae9-f668048c0c92_fila_4.png"
$ws.Range("B16").Value = "Yadder Fernando Torres"
$ws.Range("A17").Value = "db3cc54e-26a2-4d73-b2be-2b3d23ca4f1c_fila_2.png"
$ws.Range("B17").Value = "Isabella Dompe Estrada"
$ws.Range("A18").Value = "efcde127-403b-4b2d-b9a7-c10f543c35ed_fila_7.png"
$ws.Range("B18").Value = "Cristina Jozabed Carvajal"
$ws.Range("A19").Value = "f512c4d8-3979-45e8-8f96-3317e7d77d27_fila_10.png"
$ws.Range("B19").Value = "Abraham Silva Ampre"

Write-Output "Updated rows 5-19 on sheet: $($ws.Name)"
